# Fruta / hortaliza, semanal
# Inserts a new weekly price record at row 32 (date serial 44992),
# pushing all subsequent rows (old 32-40) down by one (new 33-41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 32, shifting rows 32:40 down to 33:41
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with the new record's data
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44992
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100101
$ws.Range("H32").Value = "Berries"
$ws.Range("I32").Value = 100101001
$ws.Range("J32").Value = "Arándano (blue)"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 3000
$ws.Range("O32").Value = 3000
$ws.Range("P32").Value = 3000
$ws.Range("Q32").Value = "$/bandeja 2 kilos"
$ws.Range("R32").Value = "Provincia de Diguillín"
$ws.Range("S32").Value = 1500
$ws.Range("T32").Value = 2
